# Agregado sistema de puntajes
# Adds a "Puntaje" (Score) / "Lenguaje" (Language) ranking table in columns E:F,
# keeping column D as a blank spacer column, matching the header style used
# by the existing PYPL/TIOBE/REDMONK headers in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold font, border, centered) from A1
# onto the new header cells D1:F1 so the new "spacer" + "Puntaje" + "Lenguaje"
# headers match the look of PYPL / TIOBE / REDMONK.
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E1").Value = "Puntaje"
$ws.Range("F1").Value = "Lenguaje"

# Ranking data: score (column E) and language (column F) for rows 2-11.
$scores = @(29, 26, 25, 21, 16, 10, 9, 9, 7, 7)
$languages = @("Python", "Java", "JavaScript", "C#", "PHP", "R", "Go", "TypeScript", "Rust", "Swift")

for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $scores[$i]
    $ws.Cells.Item($row, 6).Value = $languages[$i]
}
